# The scraper re-ran (400 -> 490 datapoints) and a handful of rows in the
# existing 0-39 block came back in a different order than before (the
# index column A stays a simple sequential counter, but the url/name/
# style/colorway/price/date columns for a few rows swapped places).
#
# Net effect observed between before/after workbooks: rows 2-41 (data
# rows 0-39) are a pure permutation of each other - same row contents,
# reordered. Column A (the running index) is untouched.
#
# before-row (0-based, relative to first data row) -> after-row mapping:
#   after[i] = before[mapping[i]]
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$firstRow = 2
$lastRow = 41
$rowCount = $lastRow - $firstRow + 1

# 1-based source row offsets (within the B:G block) for each destination
# row, in order.
$mapping = @(1,2,3,5,4,6,7,10,8,9,11,12,13,14,15,16,18,17,19,20,21,22,23,24,25,26,27,28,29,30,33,34,32,31,35,36,37,38,39,40)

$srcRange = $ws.Range("B$firstRow`:G$lastRow")
$values = $srcRange.Value()

$newValues = New-Object 'object[,]' $rowCount,6
for ($i = 0; $i -lt $rowCount; $i++) {
    $srcRow = $mapping[$i]
    for ($c = 1; $c -le 6; $c++) {
        $newValues[$i, $c - 1] = $values[$srcRow, $c]
    }
}

$destRange = $ws.Range("B$firstRow`:G$lastRow")
$destRange.Value = $newValues
